$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("model_5_7_24", 0.6974622203418928, -6.217302954795444, 0.5327721081510806, -0.08427227848148955, 0.3348202109336853, 1.38331663608551, 0.7297975420951843, 1.075778126716614),
    @("model_5_7_17", 0.7089548752024685, -5.737751787781138, 0.6375893158206252, 0.0424263715803016, 0.3221011757850647, 1.291402578353882, 0.5660759210586548, 0.9500721096992493),
    @("model_5_7_22", 0.7126771490506053, -5.985447917359096, 0.59439944163446, -0.01490460840409757, 0.3179817497730255, 1.33887779712677, 0.6335372924804688, 1.006953835487366),
    @("model_5_7_21", 0.7151745605062203, -5.935813290120324, 0.6011622273031679, -0.004817765460928936, 0.3152178227901459, 1.329364418983459, 0.6229740381240845, 0.9969460964202881),
    @("model_5_7_23", 0.7153249692381314, -5.91892469633464, 0.5897916877482965, -0.01151359812352126, 0.3150513470172882, 1.326127529144287, 0.6407344937324524, 1.003589510917664),
    @("model_5_7_18", 0.7205587080049182, -5.61311624607428, 0.6315835934542445, 0.05072301364502985, 0.3092591762542725, 1.267514228820801, 0.5754566192626953, 0.9418404698371887),
    @("model_5_7_19", 0.7237914803279565, -5.570136924563408, 0.6247981278214925, 0.05009121267545669, 0.305681437253952, 1.259276509284973, 0.5860553979873657, 0.9424672722816467),
    @("model_5_7_20", 0.7241726760507536, -5.587877956691899, 0.6084087030547871, 0.03613553397929159, 0.3052595555782318, 1.26267683506012, 0.6116552352905273, 0.9563136100769043),
    @("model_5_7_16", 0.7311611151266564, -4.966105110955106, 0.6563395987623593, 0.1352350975887898, 0.2975254058837891, 1.143503665924072, 0.5367884039878845, 0.8579905033111572),
    @("model_5_7_8", 0.7494968543569596, -4.510270895971352, 0.8396075823216751, 0.317627648366365, 0.2772331535816193, 1.056135416030884, 0.2505286931991577, 0.6770266890525818),
    @("model_5_7_6", 0.7502496160343931, -4.481272550906066, 0.8703661557661082, 0.3433810732465784, 0.276400089263916, 1.050577402114868, 0.2024846225976944, 0.6514750719070435),
    @("model_5_7_7", 0.7512751079756441, -4.469478495608826, 0.8584806354998793, 0.3357815107734242, 0.2752651572227478, 1.048316955566406, 0.221049502491951, 0.6590151190757751),
    @("model_5_7_15", 0.7537629801268246, -4.321862004777125, 0.7201098942714019, 0.2483666437527637, 0.2725118100643158, 1.020023822784424, 0.4371809363365173, 0.7457451820373535),
    @("model_5_7_5", 0.756672133898031, -4.184476321100356, 0.8770444387269956, 0.3786829639069402, 0.2692922651767731, 0.9936915636062622, 0.1920533329248428, 0.6164497137069702),
    @("model_5_7_14", 0.7568338979323226, -4.145829242278853, 0.7216816853537329, 0.2675346686254851, 0.2691132426261902, 0.9862841367721558, 0.434725821018219, 0.7267272472381592),
    @("model_5_7_9", 0.7583784953613861, -4.127602289378686, 0.8108148113539736, 0.3354336140673088, 0.2674038112163544, 0.9827906489372253, 0.2955022156238556, 0.6593602895736694),
    @("model_5_7_10", 0.7739916680614048, -3.546450216334809, 0.7870873816403441, 0.3772900085872629, 0.2501246333122253, 0.8714030981063843, 0.3325638473033905, 0.6178317666053772),
    @("model_5_7_13", 0.7830262461403271, -2.990722261794831, 0.7165082365252546, 0.3818373446765225, 0.2401260137557983, 0.7648885250091553, 0.4428066313266754, 0.6133200526237488),
    @("model_5_7_12", 0.7846935174009531, -2.979379301277777, 0.7323688308092395, 0.3947474523242854, 0.2382808476686478, 0.7627145648002625, 0.4180328249931335, 0.6005110740661621),
    @("model_5_7_11", 0.7850846528085643, -3.00825931285572, 0.7455792133774939, 0.401580176431654, 0.2378479391336441, 0.7682498097419739, 0.3973985314369202, 0.5937318801879883),
    @("model_5_7_0", 0.868740422980479, 0.9968286111236023, 0.996378394504709, 0.9969925229831027, 0.1452656835317612, 0.0006078496808186173, 0.005656851455569267, 0.002983917016535997),
    @("model_5_7_1", 0.873371440559328, 0.9891996301004663, 0.9942087163002719, 0.9946050369768564, 0.1401405036449432, 0.002070071408525109, 0.009045831859111786, 0.005352699663490057),
    @("model_5_7_2", 0.8763806242887576, 0.9377755261822076, 0.9882689827973903, 0.9849451377634348, 0.1368102133274078, 0.01192635949701071, 0.01832353696227074, 0.01493692584335804),
    @("model_5_7_3", 0.8794261526790388, 0.9066507389151482, 0.9846461507564183, 0.9790781588684516, 0.1334397196769714, 0.01789194531738758, 0.02398230507969856, 0.02075794152915478),
    @("model_5_7_4", 0.8798382529207642, 0.8118621695379652, 0.9764038878736441, 0.9632776125493101, 0.1329836547374725, 0.03605975583195686, 0.03685650229454041, 0.03643471002578735)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $i + 2
    $row = $data[$i]
    for ($j = 0; $j -lt $row.Count; $j++) {
        $col = $j + 1
        $ws.Cells.Item($r, $col).Value2 = $row[$j]
    }
}
